# "file upload in server"
#
# The RTC issue/status sheet's Safari note (merged cell I3:V3) is updated:
# the old note described two separate problems (Safari video playback AND
# a second line about PC browser file-upload/canvas sync). The second line
# is removed and the Safari remark is reworded to call out that the local
# video does not auto-play.
#
# The cell I3 is the top-left anchor of the merged range I3:V3, so updating
# its value updates the visible merged-cell text. We also restore the
# selection to that cell/range, matching the saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "Safari : 아이폰 영상 오류 ( 내영상 자동 재생안됨.. ), 캔버스 화면 UI 오류, 녹화 불가능 "

$ws.Range("I3:V3").Select()
